$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.138.37"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.004.83"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.61"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.97"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.11"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.140"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0854"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.00"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "3.473.22"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.63"
$ws.Range("D16").Value = "3.005.19"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.02"
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("D18").Value = "52.150.37"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.44"
$ws.Range("E19").Value = "  +5.33%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.43"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.72"
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.89"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.20"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.97"
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.18"
$ws.Range("E34").Value = "  +15.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.09"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0438"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.52"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.59"
$ws.Range("E43").Value = "  +7.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.80"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").Value = "2.120.92"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("E47").Value = "  -3.41%  "
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").Value = "3.296.96"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0332"
$ws.Range("E51").Value = "  +0.16%  "
